# Fix organism naming in the Synergy reproducibility samples/runs sheet:
#  - "Candida kruesi" (misspelling) -> "Pichia kudriavzevii" (updated species name)
#  - "Pseudominas aeruginosa" (typo) -> "Pseudomonas aeruginosa" (corrected spelling)
# These corrections apply to both the free-text "Description" column (E) and the
# "Organism" column (F) of the table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
  foreach ($col in 5, 6) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()
    if ($val -ne $null) {
      $newVal = $val -replace "Candida kruesi", "Pichia kudriavzevii"
      $newVal = $newVal -replace "Pseudominas aeruginosa", "Pseudomonas aeruginosa"
      if ($newVal -ne $val) {
        $cell.Value = $newVal
      }
    }
  }
}

# Reflect the updated selection state (whole table selected) left behind in the
# workbook after the edits were made.
$ws.Range("A1:F80").Select()
